# Fruta / hortaliza, semanal
# A new weekly observation is inserted as row 59 (pushing the existing
# rows 59-173 down to 60-174); all other columns for the new row mirror
# the row that used to be at position 59 (same market / product /
# variety / quality / unit / origin), only the date and the
# volume/price figures differ.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 59, shifting rows 59:173
# down to 60:174 (dimension grows from A1:T173 to A1:T174).
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly record.
$ws.Range("A59").Value = 10
$ws.Range("B59").Value = "Vega Modelo de Temuco"
$ws.Range("C59").Value = "La Araucanía"
$ws.Range("D59").Value = 44883
$ws.Range("E59").Value = 9
$ws.Range("F59").Value = "Fruta"
$ws.Range("G59").Value = 100107
$ws.Range("H59").Value = "Otros"
$ws.Range("I59").Value = 100107002
$ws.Range("J59").Value = "Chirimoya"
$ws.Range("K59").Value = "Cultivar IV Región"
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 20
$ws.Range("N59").Value = 2800
$ws.Range("O59").Value = 2800
$ws.Range("P59").Value = 2800
$ws.Range("Q59").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R59").Value = "Provincia del Elquí"
$ws.Range("S59").Value = 2800
$ws.Range("T59").Value = 1
